{"js": "// Office.js (Word JavaScript API) script applying the CV content edits:\n//  - Rename three section headings/subtitles in the \"Experience\" block.\n//  - Rename the \"Curated Contributions\" heading/subtitles block.\n//  - Rename the \"Featured Projects\" heading.\n//  - Fill in the previously-empty \"location\" line that follows each\n//    \"Company \u2014 date range\" line with a localized place name.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// 1. Heading / subtitle text replacements (exact paragraph text match).\nconst textReplacements = {\n  \"Signature Engagements\": \"Key Projects\",\n  \"Long-running leadership roles and programs that shape strategic outcomes.\":\n    \"Major technical and organizational engagements with long-term responsibility.\",\n  \"This is a curated selection of multi-month and multi-year engagements; more major projects are available on request.\":\n    \"A curated selection of multi-month and multi-year engagements; further details available on request.\",\n  \"Curated Contributions\": \"Additional & Focused Projects\",\n  \"Part-time, flexible, or focused initiatives that complement the bigger picture.\":\n    \"Complementary or specialized projects with flexible scope.\",\n  \"Only a handful of small projects are shown below; even more collaborations are available on request.\":\n    \"Specialized, flexible, or focused initiatives that complement the broader work.\",\n  \"Featured Projects\": \"Developed Solutions\",\n};\n\n// 2. Empty \"location\" paragraph that follows each \"Company \u2014 dates\" line,\n//    keyed by the (unique) date-range text found on that preceding line.\nconst locations = {\n  \"Jan 2024 - Sep 2025\": \"Hamburg, Germany\",\n  \"Oct 2023 - Jun 2024\": \"London, United Kingdom\",\n  \"Jan 2019 - Sep 2023\": \"Frankfurt am Main, Germany\",\n  \"Jul 2015 - Dec 2018\": \"Frankfurt am Main, Germany\",\n  \"May 2018 - Present\": \"Hockenheim, Germany\",\n  \"May 2025 - Present\": \"Seevetal, Germany\",\n  \"Aug 2025\": \"Hockenheim, Germany\",\n};\n\nconst items = paragraphs.items;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n\n  if (Object.prototype.hasOwnProperty.call(textReplacements, t)) {\n    items[i].insertText(textReplacements[t], \"Replace\");\n  }\n\n  for (const marker in locations) {\n    if (t.indexOf(marker) !== -1 && i + 1 < items.length) {\n      const nextPara = items[i + 1];\n      if (nextPara.text.trim() === \"\") {\n        nextPara.insertText(locations[marker], \"Replace\");\n      }\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Word COM interop script (PowerShell-style) applying the CV content edits:\n#  - Rename three section headings/subtitles in the \"Experience\" block.\n#  - Rename the \"Curated Contributions\" heading/subtitles block.\n#  - Rename the \"Featured Projects\" heading.\n#  - Fill in the previously-empty \"location\" line that follows each\n#    \"Company \u2014 date range\" line with a localized place name.\n#\n# NOTE: text is assigned straight onto each Paragraph's Range (found by\n# walking $d.Paragraphs and matching on its current text) rather than via\n# Find.Execute()'s Replacement, which keeps each run's existing formatting\n# (rPr) and xml:space=\"preserve\" serialization intact.\n\n$d = $word.ActiveDocument\n\n# --- 1. Heading / subtitle text replacements (exact paragraph text) -------\n$textReplacements = @(\n  @{ Find = \"Signature Engagements\"; Replace = \"Key Projects\" },\n  @{ Find = \"Long-running leadership roles and programs that shape strategic outcomes.\"; Replace = \"Major technical and organizational engagements with long-term responsibility.\" },\n  @{ Find = \"This is a curated selection of multi-month and multi-year engagements; more major projects are available on request.\"; Replace = \"A curated selection of multi-month and multi-year engagements; further details available on request.\" },\n  @{ Find = \"Curated Contributions\"; Replace = \"Additional & Focused Projects\" },\n  @{ Find = \"Part-time, flexible, or focused initiatives that complement the bigger picture.\"; Replace = \"Complementary or specialized projects with flexible scope.\" },\n  @{ Find = \"Only a handful of small projects are shown below; even more collaborations are available on request.\"; Replace = \"Specialized, flexible, or focused initiatives that complement the broader work.\" },\n  @{ Find = \"Featured Projects\"; Replace = \"Developed Solutions\" }\n)\n\n# --- 2. Empty \"location\" paragraph that follows each \"Company \u2014 dates\" ----\n#        paragraph, keyed by the (unique) date range text on that line.     -\n$locations = @{\n  \"Jan 2024 - Sep 2025\" = \"Hamburg, Germany\";\n  \"Oct 2023 - Jun 2024\" = \"London, United Kingdom\";\n  \"Jan 2019 - Sep 2023\" = \"Frankfurt am Main, Germany\";\n  \"Jul 2015 - Dec 2018\" = \"Frankfurt am Main, Germany\";\n  \"May 2018 - Present\"  = \"Hockenheim, Germany\";\n  \"May 2025 - Present\"  = \"Seevetal, Germany\";\n  \"Aug 2025\"            = \"Hockenheim, Germany\"\n}\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n  $p = $d.Paragraphs($i)\n  # Paragraph.Range.Text includes the trailing paragraph mark (chr 13) -\n  # strip it before comparing against plain strings.\n  $t = $p.Range.Text.TrimEnd([char]13)\n\n  foreach ($item in $textReplacements) {\n    if ($t -eq $item.Find) {\n      $p.Range.Text = $item.Replace\n    }\n  }\n\n  foreach ($marker in $locations.Keys) {\n    if ($t -like \"*$marker*\") {\n      $locPara = $d.Paragraphs($i + 1)\n      $locText = $locPara.Range.Text.TrimEnd([char]13)\n      if ($locText.Trim() -eq \"\") {\n        $locPara.Range.Text = $locations[$marker]\n      }\n    }\n  }\n}\n\nWrite-Output \"edit complete\"\n"}
